# Apply the edits described by the diff:
# - Update the date line
# - Update a series of "a÷b=" division problems throughout the table
# - Row containing 29÷3=,76÷7=,28÷3=,79÷6=,54÷9= is restructured into
#   19÷6=,54÷9=,18÷9=,23÷4=,40÷4=

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Date line
Replace-Text "2025-04-07 Monday" "2025-04-08 Tuesday"

# Row 1
Replace-Text "47÷4=" "26÷9="
Replace-Text "66÷4=" "86÷6="
Replace-Text "57÷5=" "60÷4="
Replace-Text "23÷5=" "91÷5="
Replace-Text "94÷2=" "90÷2="

# Row 2
Replace-Text "26÷5=" "14÷9="
Replace-Text "82÷7=" "20÷2="
Replace-Text "95÷6=" "78÷3="
Replace-Text "82÷9=" "32÷4="
Replace-Text "20÷9=" "44÷8="

# Row 3
Replace-Text "50÷2=" "86÷3="
Replace-Text "48÷9=" "13÷6="
Replace-Text "27÷3=" "24÷4="
Replace-Text "96÷2=" "70÷8="
Replace-Text "20÷5=" "49÷6="

# Row 4: cells are restructured (3 cells removed, 3 cells added), but the
# cell count in the row stays at 5. Set each cell's text directly using the
# table/cell object model to reflect the new contents in order.
$table = $d.Tables.Item(1)
$row4 = 13   # 1-indexed Word table row corresponding to this data row
$table.Cell($row4, 1).Range.Text = "19÷6="
$table.Cell($row4, 2).Range.Text = "54÷9="
$table.Cell($row4, 3).Range.Text = "18÷9="
$table.Cell($row4, 4).Range.Text = "23÷4="
$table.Cell($row4, 5).Range.Text = "40÷4="

# Row 5
Replace-Text "96÷5=" "81÷9="
Replace-Text "83÷6=" "12÷6="
Replace-Text "60÷5=" "42÷7="
Replace-Text "75÷4=" "45÷3="
Replace-Text "39÷2=" "19÷2="
